$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-09-30 00:00:00"
$ws.Range("O2").Value = 26941818.92
$ws.Range("P2").Value = 883.2047933684
$ws.Range("Q2").Value = 136864952.75
$ws.Range("R2").Value = 4486.6971555215
$ws.Range("S2").Value = 46256312.16
$ws.Range("T2").Value = 1516.3711382875
$ws.Range("U2").Value = 731181.4
$ws.Range("V2").Value = 23.9695366976
$ws.Range("W2").Value = 1144324.17
$ws.Range("X2").Value = 37.5131536261
$ws.Range("Y2").Value = 2813142.77
$ws.Range("Z2").Value = 92.2202463862
$ws.Range("AA2").Value = -24606472.93
$ws.Range("AB2").Value = -806.6476470722999
$ws.Range("AC2").Value = 3050461.13
$ws.Range("AD2").Value = ""
